$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (day-of-month labels) for columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update meanEMG row (row 2) for columns B:E
$ws.Range("B2").Value = -2.2480288825360719
$ws.Range("C2").Value = 3.3673641978477207
$ws.Range("D2").Value = 5.4829189350170298
$ws.Range("E2").Value = 10.152546027031246

# Update legmaxROM row (row 3) for columns B:E
$ws.Range("B3").Value = -7.2455310004589952
$ws.Range("C3").Value = 4.6464066055564892
$ws.Range("D3").Value = 12.721568865247699
$ws.Range("E3").Value = 5.0178586006840877

# Update the saved selection to match the new, narrower range of interest
$ws.Activate()
$ws.Range("B1:E3").Select()
